$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G (shifts old G:K -> H:L), mirroring the
# worksheet's original "Insert Column" behaviour (new col G inherits F's
# formatting by default).
$ws.Columns("G").EntireColumn.Insert()

# Make sure the new column G exactly matches column F's cell formatting
# (style) for the header rows, same as Excel's native insert-column copy.
$ws.Range("F2:F5").Copy()
$ws.Range("G2:G5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New column G content: "Archival Object URI" field, paralleling the
# existing "Ref_id" (ao_ref_id) field in column F.
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = "Archival Object URI  REQUIRED IF NO REF ID"
$ws.Range("G4").Value = "ao_uri"
$ws.Range("G5").Value = "Arch. Obj. URI"

# Column F's description (row 3) now clarifies it is only required when
# no Archival Object URI is supplied.
$ws.Range("F3").Value = "Ref_id  REQUIRED IF NO URI"

# Match column G's width to column F's.
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Reflect the final selection/active cell left on the sheet.
$ws.Range("G4").Select()
